{"js": "// Replace the two-digit-division answers in the document's table, in\n// document order. The same \"old\" text can appear more than once (and some\n// \"new\" values repeat too), so each pair is applied to the first remaining\n// matching occurrence \u2014 never touching cells already updated by an earlier\n// pair in this list.\nconst replacements = [\n  [\"13\u00f75=2, 3\", \"71\u00f76=11, 5\"],\n  [\"16\u00f79=1, 7\", \"68\u00f73=22, 2\"],\n  [\"66\u00f77=9, 3\", \"67\u00f74=16, 3\"],\n  [\"64\u00f78=8, 0\", \"61\u00f75=12, 1\"],\n  [\"89\u00f76=14, 5\", \"77\u00f77=11, 0\"],\n  [\"54\u00f79=6, 0\", \"71\u00f76=11, 5\"],\n  [\"52\u00f76=8, 4\", \"59\u00f73=19, 2\"],\n  [\"27\u00f74=6, 3\", \"66\u00f77=9, 3\"],\n  [\"51\u00f77=7, 2\", \"67\u00f75=13, 2\"],\n  [\"30\u00f73=10, 0\", \"25\u00f79=2, 7\"],\n  [\"16\u00f78=2, 0\", \"91\u00f78=11, 3\"],\n  [\"78\u00f74=19, 2\", \"31\u00f75=6, 1\"],\n  [\"89\u00f79=9, 8\", \"94\u00f76=15, 4\"],\n  [\"29\u00f74=7, 1\", \"79\u00f74=19, 3\"],\n  [\"32\u00f74=8, 0\", \"58\u00f74=14, 2\"],\n  [\"25\u00f74=6, 1\", \"69\u00f72=34, 1\"],\n  [\"25\u00f74=6, 1\", \"94\u00f79=10, 4\"],\n  [\"86\u00f77=12, 2\", \"75\u00f74=18, 3\"],\n  [\"58\u00f72=29, 0\", \"94\u00f76=15, 4\"],\n  [\"15\u00f77=2, 1\", \"71\u00f72=35, 1\"],\n  [\"53\u00f73=17, 2\", \"42\u00f77=6, 0\"],\n  [\"31\u00f79=3, 4\", \"44\u00f79=4, 8\"],\n  [\"53\u00f79=5, 8\", \"19\u00f72=9, 1\"],\n  [\"90\u00f73=30, 0\", \"10\u00f76=1, 4\"],\n  [\"21\u00f77=3, 0\", \"26\u00f77=3, 5\"],\n];\n\n// Each iteration re-searches the (already partially edited) document for\n// the next pair's \"old\" text and replaces the FIRST remaining hit. Because\n// earlier pairs in the list are always applied first, any previously\n// replaced occurrence no longer matches \"old\" text, so a fresh items[0] is\n// always the correct (next, left-to-right) occurrence \u2014 this naturally\n// handles the repeated \"25\u00f74=6, 1\" old text without needing an explicit\n// \"skip N already-used hits\" counter.\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    // Nothing left to replace for this text \u2014 skip defensively.\n    continue;\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-division answers in the document's table, in\n# document order. The same \"old\" text can appear more than once (and some\n# \"new\" values repeat too), so each pair is applied with wdReplaceOne\n# (replace a single occurrence) against a fresh Content range each time \u2014\n# this always lands on the first REMAINING match, since earlier pairs in\n# this list are applied first and therefore already consumed the earlier\n# occurrences.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"13\u00f75=2, 3\", \"71\u00f76=11, 5\"),\n    @(\"16\u00f79=1, 7\", \"68\u00f73=22, 2\"),\n    @(\"66\u00f77=9, 3\", \"67\u00f74=16, 3\"),\n    @(\"64\u00f78=8, 0\", \"61\u00f75=12, 1\"),\n    @(\"89\u00f76=14, 5\", \"77\u00f77=11, 0\"),\n    @(\"54\u00f79=6, 0\", \"71\u00f76=11, 5\"),\n    @(\"52\u00f76=8, 4\", \"59\u00f73=19, 2\"),\n    @(\"27\u00f74=6, 3\", \"66\u00f77=9, 3\"),\n    @(\"51\u00f77=7, 2\", \"67\u00f75=13, 2\"),\n    @(\"30\u00f73=10, 0\", \"25\u00f79=2, 7\"),\n    @(\"16\u00f78=2, 0\", \"91\u00f78=11, 3\"),\n    @(\"78\u00f74=19, 2\", \"31\u00f75=6, 1\"),\n    @(\"89\u00f79=9, 8\", \"94\u00f76=15, 4\"),\n    @(\"29\u00f74=7, 1\", \"79\u00f74=19, 3\"),\n    @(\"32\u00f74=8, 0\", \"58\u00f74=14, 2\"),\n    @(\"25\u00f74=6, 1\", \"69\u00f72=34, 1\"),\n    @(\"25\u00f74=6, 1\", \"94\u00f79=10, 4\"),\n    @(\"86\u00f77=12, 2\", \"75\u00f74=18, 3\"),\n    @(\"58\u00f72=29, 0\", \"94\u00f76=15, 4\"),\n    @(\"15\u00f77=2, 1\", \"71\u00f72=35, 1\"),\n    @(\"53\u00f73=17, 2\", \"42\u00f77=6, 0\"),\n    @(\"31\u00f79=3, 4\", \"44\u00f79=4, 8\"),\n    @(\"53\u00f79=5, 8\", \"19\u00f72=9, 1\"),\n    @(\"90\u00f73=30, 0\", \"10\u00f76=1, 4\"),\n    @(\"21\u00f77=3, 0\", \"26\u00f77=3, 5\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue = 1, wdReplaceOne = 1\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null\n}\n"}
